# Generate Report for Handback
# Update the handoff/handback timestamps on the zh-cn and de-de sheets
# to reflect a freshly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 04:54:06"
$wsZhCn.Range("H2").Value = "2016-03-20 04:54:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 04:54:09"
$wsDeDe.Range("H2").Value = "2016-03-20 04:54:30"
